# Automatische test-sync: 2025-06-19 22:29:50
# Appends two new mail-log rows to the "Logs" sheet and bumps the matching
# category counters on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Row 41: Productinformatie -------------------------------------------
$logs.Cells.Item(41, 1).Value = "Is product X op voorraad?"
$logs.Cells.Item(41, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(41, 3).Value = "Ik ben geïnteresseerd in product X. Is dit momenteel op voorraad?"
$logs.Cells.Item(41, 4).Value = "Productinformatie"
$logs.Cells.Item(41, 6).Value = "2025-06-19 22:29:11"
$logs.Cells.Item(41, 7).Value = "Nee"

# --- Row 42: Sollicitatie / Vacature --------------------------------------
$logs.Cells.Item(42, 1).Value = "Sollicitatie marketingfunctie"
$logs.Cells.Item(42, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(42, 3).Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$logs.Cells.Item(42, 4).Value = "Sollicitatie / Vacature"
$logs.Cells.Item(42, 6).Value = "2025-06-19 22:29:12"
$logs.Cells.Item(42, 7).Value = "Nee"

# --- Dashboard counters ----------------------------------------------------
# Productinformatie: 5 -> 6
$dash.Cells.Item(3, 2).Value = 6
# Sollicitatie / Vacature: 1 -> 2
$dash.Cells.Item(10, 2).Value = 2

# --- Extend conditional-formatting ranges to cover the new rows -----------
$catRules = $logs.Range("D2:D40").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D42"))
}

$answeredRules = $logs.Range("G2:G40").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G42"))
}
